$d = $word.ActiveDocument

# Locate the two target paragraphs by their text ("Package usercmd" / "Package scontrol")
$pUserCmd = $null
$pScontrol = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -eq "Package usercmd`r") {
        $pUserCmd = $i
    } elseif ($txt -eq "Package scontrol`r") {
        $pScontrol = $i
    }
}

if (-not $pUserCmd) { throw "Could not find paragraph 'Package usercmd'" }
if (-not $pScontrol) { throw "Could not find paragraph 'Package scontrol'" }

# "Package usercmd" paragraph -> becomes a Titolo2 heading, and both its runs get
# light-gray highlighting.
$pu = $d.Paragraphs.Item($pUserCmd)
$pu.Style = "Titolo2"
$pu.Range.HighlightColorIndex = 16

# Insert two new plain paragraphs ("Interruttore", "Timer") right after "Package scontrol"
# *before* restyling "Package scontrol", so the new paragraphs don't inherit Titolo2.
$ps = $d.Paragraphs.Item($pScontrol)
$r = $ps.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$pInterruttore = $d.Paragraphs.Item($pScontrol + 1)
$pInterruttore.Range.Text = "Interruttore"

$r2 = $pInterruttore.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$pTimer = $d.Paragraphs.Item($pScontrol + 2)
$pTimer.Range.Text = "Timer"

# "Package scontrol" paragraph -> becomes a Titolo2 heading (no highlight change).
$ps.Style = "Titolo2"
